$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.775.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.50%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.568.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.03%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.98%  "

# Row 6
$ws.Range("E6").Value = "  -2.05%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.95"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.00%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.248"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.86%  "

# Row 10
$ws.Range("E10").Value = "  -1.37%  "

# Row 11
$ws.Range("E11").Value = "  -0.39%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.790.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.03%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.584.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.25%  "

# Row 14
$ws.Range("E14").Value = "  -2.57%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.515"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.87%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.790.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.46%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.61%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.80%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.51%  "

# Row 20
$ws.Range("E20").Value = "  -1.91%  "

# Row 21
$ws.Range("E21").Value = "  +0.07%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.03%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.54%  "

# Row 24
$ws.Range("E24").Value = "  -0.92%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.16%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.66%  "

# Row 27
$ws.Range("E27").Value = "  -0.34%  "

# Row 28
$ws.Range("E28").Value = "  +0.06%  "

# Row 29
$ws.Range("E29").Value = "  -1.52%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0465"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.06%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.02%  "

# Row 32
$ws.Range("E32").Value = "  -1.38%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.395.52"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.35%  "

# Row 34
$ws.Range("E34").Value = "  -1.42%  "

# Row 35
$ws.Range("E35").Value = "  -0.60%  "

# Row 36
$ws.Range("E36").Value = "  -1.17%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.932"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.72%  "

# Row 38
$ws.Range("E38").Value = "  -3.02%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.528"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.45%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.818"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.11%  "

# Row 41
$ws.Range("E41").Value = "  +0.10%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.989"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.00%  "

# Row 43
$ws.Range("E43").Value = "  -0.31%  "

# Row 44
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.19"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.21%  "

# Row 45
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.08%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.52%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.703.70"
$ws.Range("D47").Style = "Normal"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.87%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₇0983"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.48%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0954"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.42%  "

# Row 51
$ws.Range("E51").Value = "  -0.96%  "
